$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers (now human-readable Spanish labels, reordered)
$ws.Cells.Item(1,1).Value = "Número de edificios"
$ws.Cells.Item(1,2).Value = "Comarca nombre"
$ws.Cells.Item(1,3).Value = "Comarca código"
$ws.Cells.Item(1,4).Value = "Provincia código"
$ws.Cells.Item(1,5).Value = "Aragón"
$ws.Cells.Item(1,6).Value = "Municipio código"
$ws.Cells.Item(1,7).Value = "Provincia nombre"
$ws.Cells.Item(1,8).Value = "Año de construcción"
$ws.Cells.Item(1,9).Value = "Municipio nombre"

# Row 2: sdmx/iaest dimension metadata
$ws.Cells.Item(2,1).Value = "iaest-measure:numero-de-edificios"
$ws.Cells.Item(2,2).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2,3).Value = "null"
$ws.Cells.Item(2,4).Value = "null"
$ws.Cells.Item(2,5).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2,6).Value = "null"
$ws.Cells.Item(2,7).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2,8).Value = "iaest-dimension:ano-de-construccion"
$ws.Cells.Item(2,9).Value = "sdmx-dimension:refArea"

# Row 3: dim/medida classifiers
$ws.Cells.Item(3,1).Value = "medida"
$ws.Cells.Item(3,2).Value = "dim"
$ws.Cells.Item(3,3).Value = "null"
$ws.Cells.Item(3,4).Value = "null"
$ws.Cells.Item(3,5).Value = "dim"
$ws.Cells.Item(3,6).Value = "null"
$ws.Cells.Item(3,7).Value = "dim"
$ws.Cells.Item(3,8).Value = "dim"
$ws.Cells.Item(3,9).Value = "dim"

# Row 4: data types / URI references
$ws.Cells.Item(4,1).Value = "xsd:int"
$ws.Cells.Item(4,2).Value = "URI-comarca"
$ws.Cells.Item(4,3).Value = "null"
$ws.Cells.Item(4,4).Value = "null"
$ws.Cells.Item(4,5).Value = "URI-Comunidad"
$ws.Cells.Item(4,6).Value = "null"
$ws.Cells.Item(4,7).Value = "URI-Provincia"
$ws.Cells.Item(4,8).Value = "skos:Concept"
$ws.Cells.Item(4,9).Value = "URI-Municipio"

# Row 5: the mapping file reference moves from A5 to H5
# Copy formatting from an already-styled cell so H5 keeps the same cell style as its neighbours
$ws.Cells.Item(4,8).Copy()
$ws.Cells.Item(5,8).PasteSpecial(-4122)
$ws.Cells.Item(5,8).Value = "mapping-ano-de-construccion.xlsx"
$ws.Cells.Item(5,1).ClearContents()
